$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.40%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'35.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'13.02%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.169"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.92%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07794"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'2.388"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'8.59%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.057"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.42%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.969"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'6.27%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9301"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.03%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09985"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'8.72%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'9.96%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08672"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03311"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'6.28%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09902"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.51%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001484"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.91%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005783"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.49%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-1.60%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.130"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.52%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'1.35%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1325"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.49%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.323"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'4.06%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'13.57%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04566"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.64%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.49%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004447"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.95%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.18%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003695"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'8.80%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01782"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'13.68%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04800"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'8.12%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007746"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.58%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'6.51%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007158"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-20.35%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002186"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.01%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'4.81%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006100"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.40%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'17.86%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = "Normal"
